$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ingresos")

$ws.Range("A1").Value = "Fecha"
$ws.Range("B1").Value = "Tipo de Pago"
$ws.Range("C1").Value = "Cantidad ($)"
$ws.Range("D1").Value = "Tipo de Pago"
$ws.Range("E1").Value = "Descripcion"
